$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Fun times"
[void]$ws.Range("A7").Select()
